# Benchmark update: 2025-12-23 06:43:16 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ŞANS OYUNLARI) - DENIZBANK (J) now has a value
$ws.Range("J2").Value = "25 TL - 25 TL"

# Row 3 (HESAPTAN EFT - Şube) - İŞBANKASI (E) cleared
$ws.Range("E3").Value = ""

# Row 4 (HESAPTAN EFT - ATM) - İŞBANKASI (E) cleared
$ws.Range("E4").Value = ""

# Row 5 (HESAPTAN EFT - Mobil) - İŞBANKASI (E) cleared
$ws.Range("E5").Value = ""

# Row 6 (DÜZENLİ EFT) - İŞBANKASI (E) cleared
$ws.Range("E6").Value = ""

# Row 7 (KREDİ KARTINDAN FATURA ÖDEME) - DENIZBANK (J) now has a value
$ws.Range("J7").Value = "%2,5"

# Row 8 (HESAPTAN HAVALE - Şube) - İŞBANKASI (E) cleared
$ws.Range("E8").Value = ""

# Row 9 (HESAPTAN HAVALE - ATM) - İŞBANKASI (E) cleared
$ws.Range("E9").Value = ""

# Row 10 (HESAPTAN HAVALE - Mobil) - İŞBANKASI (E) cleared
$ws.Range("E10").Value = ""

# Row 11 (DÜZENLİ HAVALE) - İŞBANKASI (E) cleared
$ws.Range("E11").Value = ""

# Row 13 (GELEN SWIFT)
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("E13").Value = ""
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

# Row 14 (GİDEN SWIFT - Mobil) - İŞBANKASI (E) cleared
$ws.Range("E14").Value = ""
